$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a hyperlink on A2 pointing to the Audible listener page, with the
# link URL itself as the display text (Excel's default when
# TextToDisplay is not otherwise customized). Do this before renaming the
# cell's text so the hyperlink's stored display text is the URL.
$audibleUrl = "https://www.audible.co.uk/listener/A89XFJ7UQTS9L?pf_rd_p=d85bb0c7-d987-483d-acca-afcf5c6bc241&pf_rd_r=7WM9F89SZSBNP94R9B4V&ref=a_pd_Anger-_c16_rvlsnl_0"
$ws.Hyperlinks.Add($ws.Range("A2"), $audibleUrl, "", "", $audibleUrl)

# Adding the hyperlink auto-applies the built-in "Hyperlink" cell style;
# restore the cell's original (unstyled) formatting.
$ws.Range("A2").Style = "Normal"

# Rename reviewer "Chester robin" -> "Mcintosh" (A2's visible text)
$ws.Range("A2").Value = "Mcintosh"

# Update the review date in B2 (date serial 43749 -> 43765)
$ws.Range("B2").Value = 43765

# Move the active selection from B2 to C3
$ws.Range("C3").Select() | Out-Null

# Switch the sheet's print orientation to portrait
$ws.PageSetup.Orientation = 1 | Out-Null
